$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1335.3334
$ws.Range("I41").Value = 1592.4
$ws.Range("J41").Value = 50
$ws.Range("K41").Value = 1592.4
$ws.Range("L41").Value = 50
$ws.Range("M41").Value = -1152.4
$ws.Range("N41").Value = -930
$ws.Range("H62").Value = 17150.85
$ws.Range("I62").Value = 23002.363
$ws.Range("K62").Value = 23002.363
$ws.Range("M62").Value = -22378.363
$ws.Range("H65").Value = 17150.85
$ws.Range("I65").Value = 23002.363
$ws.Range("K65").Value = 115011.815
$ws.Range("M65").Value = -111891.815
$ws.Range("H76").Value = 5745.9
$ws.Range("I76").Value = 5065.7144
$ws.Range("K76").Value = 5065.7144
$ws.Range("M76").Value = -4750.7144
$ws.Range("H79").Value = 5745.9
$ws.Range("I79").Value = 5065.7144
$ws.Range("K79").Value = 5065.7144
$ws.Range("M79").Value = -3973.7144
$ws.Range("H86").Value = 3100.6667
$ws.Range("I86").Value = 2929.3809
$ws.Range("K86").Value = 2929.3809
$ws.Range("M86").Value = -1806.3809
$ws.Range("H89").Value = 3100.6667
$ws.Range("I89").Value = 2929.3809
$ws.Range("K89").Value = 14646.9045
$ws.Range("M89").Value = -9030.904500000001
$ws.Range("H103").Value = 459.91666
$ws.Range("I103").Value = 290
$ws.Range("J103").Value = 516.55554
$ws.Range("K103").Value = 870
$ws.Range("L103").Value = 1549.66662
$ws.Range("M103").Value = -284
$ws.Range("N103").Value = -2721.66662
$ws.Range("H106").Value = 3771.818
$ws.Range("I106").Value = 3771.818
$ws.Range("K106").Value = 3771.818
$ws.Range("M106").Value = -3140.818
$ws.Range("H107").Value = 1306.7778
$ws.Range("I107").Value = 838.43475
$ws.Range("K107").Value = 838.43475
$ws.Range("M107").Value = 1081.56525
$ws.Range("H128").Value = 69749.75
$ws.Range("J128").Value = 98999
$ws.Range("L128").Value = 98999
$ws.Range("N128").Value = -108959
$ws.Range("H135").Value = 1781.1666
$ws.Range("I135").Value = 1166.3572
$ws.Range("J135").Value = 3933
$ws.Range("K135").Value = 10497.2148
$ws.Range("L135").Value = 35397
$ws.Range("M135").Value = -7962.2148
$ws.Range("N135").Value = -40467
$ws.Range("H137").Value = 2022.6
$ws.Range("I137").Value = 1230
$ws.Range("K137").Value = 3690
$ws.Range("M137").Value = -1140
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4006.6206
$ws.Range("I122").Value = 3872.9614
$ws.Range("K122").Value = 11618.8842
$ws.Range("M122").Value = -9168.8842
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 26282.25
$ws.Range("J21").Value = 26282.25
$ws.Range("L21").Value = 26282.25
$ws.Range("N21").Value = -26754.25
$ws.Range("H134").Value = 2016.7333
$ws.Range("I134").Value = 2071.0232
$ws.Range("K134").Value = 6213.069600000001
$ws.Range("M134").Value = -3678.069600000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 517.41174
$ws.Range("I22").Value = 322.36365
$ws.Range("K22").Value = 322.36365
$ws.Range("M22").Value = 27.63634999999999
$ws.Range("H31").Value = 2598.6511
$ws.Range("I31").Value = 1535.6857
$ws.Range("J31").Value = 7249.125
$ws.Range("K31").Value = 1535.6857
$ws.Range("L31").Value = 7249.125
$ws.Range("M31").Value = -1240.6857
$ws.Range("N31").Value = -7839.125
$ws.Range("H34").Value = 2598.6511
$ws.Range("I34").Value = 1535.6857
$ws.Range("J34").Value = 7249.125
$ws.Range("K34").Value = 1535.6857
$ws.Range("L34").Value = 7249.125
$ws.Range("M34").Value = -1333.6857
$ws.Range("N34").Value = -7653.125
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 49996.332
$ws.Range("J39").Value = 49996.332
$ws.Range("L39").Value = 49996.332
$ws.Range("N39").Value = -51060.332
$ws.Range("H98").Value = 15669.833
$ws.Range("J98").Value = 15669.833
$ws.Range("L98").Value = 15669.833
$ws.Range("N98").Value = -21659.833
$ws.Range("H105").Value = 75999.2
$ws.Range("J105").Value = 75999.2
$ws.Range("L105").Value = 75999.2
$ws.Range("N105").Value = -82987.2
$ws.Range("H113").Value = 143990.28
$ws.Range("J113").Value = 501499.5
$ws.Range("L113").Value = 501499.5
$ws.Range("N113").Value = -505839.5
$ws.Range("H132").Value = 2149.087
$ws.Range("I132").Value = 1216.2285
$ws.Range("K132").Value = 3648.6855
$ws.Range("M132").Value = -1118.6855
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7115.8887
$ws.Range("I7").Value = 7425.8667
$ws.Range("J7").Value = 5566
$ws.Range("K7").Value = 7425.8667
$ws.Range("L7").Value = 5566
$ws.Range("M7").Value = -7313.8667
$ws.Range("N7").Value = -5790
$ws.Range("H46").Value = 8724.538
$ws.Range("J46").Value = 4499
$ws.Range("L46").Value = 4499
$ws.Range("N46").Value = -4875
$ws.Range("H68").Value = 4875
$ws.Range("I68").Value = 2175
$ws.Range("K68").Value = 2175
$ws.Range("M68").Value = -1426
$ws.Range("H71").Value = 4875
$ws.Range("I71").Value = 2175
$ws.Range("K71").Value = 10875
$ws.Range("M71").Value = -7131
$ws.Range("H82").Value = 1990.8462
$ws.Range("I82").Value = 978.7059
$ws.Range("J82").Value = 2772.9546
$ws.Range("K82").Value = 978.7059
$ws.Range("L82").Value = 2772.9546
$ws.Range("M82").Value = -617.7059
$ws.Range("N82").Value = -3494.9546
$ws.Range("H85").Value = 1990.8462
$ws.Range("I85").Value = 978.7059
$ws.Range("J85").Value = 2772.9546
$ws.Range("K85").Value = 978.7059
$ws.Range("L85").Value = 2772.9546
$ws.Range("M85").Value = 269.2941
$ws.Range("N85").Value = -5268.9546
$ws.Range("H122").Value = 4323.241
$ws.Range("I122").Value = 3497.2307
$ws.Range("K122").Value = 10491.6921
$ws.Range("M122").Value = -8041.6921
$ws.Range("H126").Value = 7115.8887
$ws.Range("I126").Value = 7425.8667
$ws.Range("J126").Value = 5566
$ws.Range("K126").Value = 22277.6001
$ws.Range("L126").Value = 16698
$ws.Range("M126").Value = -19807.6001
$ws.Range("N126").Value = -21638
$ws.Range("H136").Value = 3163.7693
$ws.Range("I136").Value = 3107.2727
$ws.Range("K136").Value = 9321.8181
$ws.Range("M136").Value = -6771.8181
$ws.Range("H140").Value = 89999.5
$ws.Range("J140").Value = 89999.5
$ws.Range("L140").Value = 89999.5
$ws.Range("N140").Value = -100359.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 12299.8
$ws.Range("I4").Value = 29999.5
$ws.Range("K4").Value = 29999.5
$ws.Range("M4").Value = -29886.5
